$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.281.29"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.358.64"
$ws.Range("E3").Value = "  +5.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.25%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.72"
$ws.Range("E5").Value = "  +1.37%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.03%  "

# Row 7 - Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.36"
$ws.Range("E7").Value = "  +13.51%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  +13.64%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("E10").Value = "  +2.67%  "

# Row 11 - Avalanche
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.36"
$ws.Range("E11").Value = "  -0.97%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.712.08"
$ws.Range("E12").Value = "  +5.59%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.17%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.28"
$ws.Range("E14").Value = "  +5.63%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +3.84%  "

# Row 16 - Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.866"
$ws.Range("E16").Value = "  +4.92%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.358.49"
$ws.Range("E17").Value = "  +5.49%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.252.70"
$ws.Range("E18").Value = "  +0.23%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +3.78%  "

# Row 20 - was Uniswap, now Litecoin
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "74.88"
$ws.Range("E20").Value = "  +2.89%  "

# Row 21 - was Litecoin, now Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.36"
$ws.Range("E21").Value = "  +4.59%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.15"
$ws.Range("E22").Value = "  +2.15%  "

# Row 23 - WEMIXToken
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.87"
$ws.Range("E23").Value = "  +5.74%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.07%  "

# Row 25 - PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +1.86%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  +2.53%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +3.21%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.52"
$ws.Range("E28").Value = "  +4.81%  "

# Row 29 - Monero
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.75"
$ws.Range("E29").Value = "  -0.21%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +10.24%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +3.71%  "

# Row 32 - Stellar
$ws.Range("E32").Value = "  +3.02%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.79%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +2.49%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("E35").Value = "  +3.57%  "

# Row 36 - RenderToken
$ws.Range("E36").Value = "  +4.12%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +7.32%  "

# Row 38 - THORChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.52"
$ws.Range("E38").Value = "  +3.93%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +2.09%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.05"
$ws.Range("E40").Value = "  +12.34%  "

# Row 41 - BinanceUSD
$ws.Range("E41").Value = "  +0.14%  "

# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.92"
$ws.Range("E42").Value = "  +3.34%  "

# Row 43 - Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.34"
$ws.Range("E43").Value = "  +3.20%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  +9.62%  "

# Row 45 - FTXToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.49"
$ws.Range("E45").Value = "  +0.77%  "

# Row 46 - Cronos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0958"
$ws.Range("E46").Value = "  +1.96%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +2.51%  "

# Row 48 - Maker
$ws.Range("D48").Value = "1.444.32"
$ws.Range("E48").Value = "  +0.17%  "

# Row 49 - RocketPoolETH
$ws.Range("E49").Value = "  +5.76%  "

# Row 50 - HuobiToken
$ws.Range("E50").Value = "  +1.15%  "

# Row 51 - TerraClassic
$ws.Range("E51").Value = "  -3.71%  "
